$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 156, pushing existing row 156 (and below) down to 157+
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new record
# (same market/category info as the adjacent record, new date/volume/unit)
$ws.Range("A156").Value = 5
$ws.Range("B156").Value = "Macroferia Regional de Talca"
$ws.Range("C156").Value = "Maule"
$ws.Range("D156").Value = 44673
$ws.Range("E156").Value = 7
$ws.Range("F156").Value = 100112003
$ws.Range("G156").Value = "Ajo"
$ws.Range("H156").Value = "Chino"
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 300
$ws.Range("K156").Value = 20000
$ws.Range("L156").Value = 20000
$ws.Range("M156").Value = 20000
$ws.Range("N156").Value = "`$/caja 10 kilos"
$ws.Range("O156").Value = "China"
$ws.Range("P156").Value = 2000
$ws.Range("Q156").Value = 10
$ws.Range("R156").Value = "Hortaliza"
